# Adds two new category rows ("criminal" / "incidents", under a new
# "general" parent) to the COMMON sheet of the categories workbook so
# admin crime/incident notifications resolve to a known category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COMMON")

# --- Row 98: criminal -------------------------------------------------
$ws.Rows.Item(98).Insert()
$ws.Cells.Item(98, 1).Value = "criminal"
$ws.Cells.Item(98, 2).Value = "general"
$ws.Cells.Item(98, 3).Value = "Криминал"
$ws.Cells.Item(98, 4).Value = "Criminal"
$ws.Cells.Item(98, 5).Value = "Кримінал"

# --- Row 99: incidents -------------------------------------------------
$ws.Rows.Item(99).Insert()
$ws.Cells.Item(99, 1).Value = "incidents"
$ws.Cells.Item(99, 2).Value = "general"
$ws.Cells.Item(99, 3).Value = "Происшествия"
$ws.Cells.Item(99, 4).Value = "Incidents"
$ws.Cells.Item(99, 5).Value = "Інциденти"
